$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the now-unused tail rows -------------------------------------
# Old rows 22-31 (M21..M30) are gone; this shifts nothing else and shrinks
# the used range from A1:D31 down to A1:D21.
$ws.Rows("22:31").Delete()

# --- 2. Plain text labels (A/B columns for the relabelled rows) -----------
# These values never look like numbers, so a direct assignment is safe and
# leaves no stray style behind.
$ws.Cells.Item(16, 1).Value = "M20"
$ws.Cells.Item(16, 2).Value = "MMM20"
$ws.Cells.Item(17, 1).Value = "M21"
$ws.Cells.Item(17, 2).Value = "MMM21"
$ws.Cells.Item(18, 1).Value = "M22"
$ws.Cells.Item(18, 2).Value = "MMM22"
$ws.Cells.Item(19, 1).Value = "M23"
$ws.Cells.Item(19, 2).Value = "MMM23"
$ws.Cells.Item(20, 1).Value = "M24"
$ws.Cells.Item(20, 2).Value = "MMM24"
$ws.Cells.Item(21, 1).Value = "M25"
$ws.Cells.Item(21, 2).Value = "MMM25"

# --- 3. shock / extreme_level columns (C & D) ------------------------------
# Many of these look like numbers/percentages ("20.0 %", "1.1 %", ...), and
# a bare .Value assignment gets "smart" re-parsed into a percent number by
# Excel. Prefixing with an apostrophe forces literal text for every cell in
# one pass; a single bulk Style reset afterwards clears the quote-prefix
# formatting flag off the whole block again (cheaper than doing it cell by
# cell, and keeps the style table from growing once per cell).
#
# D2:D11 already hold blank text in the source workbook and the diff leaves
# them untouched, so they are intentionally skipped below.
$ws.Cells.Item(2, 3).Value  = "'20.0 %"
$ws.Cells.Item(3, 3).Value  = "'4.3 %"
$ws.Cells.Item(4, 3).Value  = "'2.4 %"
$ws.Cells.Item(5, 3).Value  = "'1.7 %"
$ws.Cells.Item(6, 3).Value  = "'1.3 %"
$ws.Cells.Item(7, 3).Value  = "'1.1 %"
$ws.Cells.Item(8, 3).Value  = "'0.9 %"
$ws.Cells.Item(9, 3).Value  = "'0.8 %"
$ws.Cells.Item(10, 3).Value = "'0.7 %"
$ws.Cells.Item(11, 3).Value = "'7.8 %"

$ws.Cells.Item(12, 3).Value = "'13% max"
$ws.Cells.Item(12, 4).Value = "'(+198 ppts)"

$ws.Cells.Item(13, 3).Value = "'13 peak"
$ws.Cells.Item(13, 4).Value = "'(+216 ppts)"

$ws.Cells.Item(14, 3).Value = "'13% peak"
$ws.Cells.Item(14, 4).Value = "'(+23400 bps)"

$ws.Cells.Item(15, 3).Value = "'13% peak"
$ws.Cells.Item(15, 4).Value = "'(+25200 bps)"

$ws.Cells.Item(16, 3).Value = "'1.1 %"
$ws.Cells.Item(16, 4).Value = "'1.2 %"

$ws.Cells.Item(17, 3).Value = "'13% peak"
$ws.Cells.Item(17, 4).Value = "'(+37800 bps)"

$ws.Cells.Item(18, 3).Value = "'13% peak"
$ws.Cells.Item(18, 4).Value = "'(+39600 bps)"

$ws.Cells.Item(19, 3).Value = "'13% peak"
$ws.Cells.Item(19, 4).Value = "'(+41400 bps)"

$ws.Cells.Item(20, 3).Value = "'0.2 %"
$ws.Cells.Item(20, 4).Value = "'"

$ws.Cells.Item(21, 3).Value = "'13% peak"
$ws.Cells.Item(21, 4).Value = "'(+45000 bps)"

# One bulk style reset over the whole edited block removes the quote-prefix
# flag Excel attached while literal-parsing the values above (D2:D11 are
# untouched in content but included here too - a Style-only touch does not
# change their blank-text value).
$ws.Range("A2:D21").Style = "Normal"
